$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with plain default style (no explicit style / quotePrefix),
# used to reset style on cells where we had to force text entry via a leading
# apostrophe (prevents Excel's auto-numeric-conversion of digit-only strings).
$plainStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "92.022.05"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "3.103.40"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'239.51"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'614.09"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  -4.73%  "
$ws.Range("D8").Value = "'0.389"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +7.25%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "3.107.09"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'0.728"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "91.919.37"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.51"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'34.13"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").Value = "3.685.88"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "3.090.82"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'3.66"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'14.72"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "'5.82"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'447.10"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'9.27"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "'11.69"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "3.275.45"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +12.73%  "
$ws.Range("D31").Value = "'0.227"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -6.98%  "
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").Value = "'9.24"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +57.31%  "
$ws.Range("D35").Value = "'0.166"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").Value = "'8.03"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("D37").Value = "'26.12"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'4.15"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("D39").Value = "'1.92"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "'482.01"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("D41").Value = "'1.30"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "'3.46"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'0.432"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'158.96"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("D47").Value = "'1.90"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'0.696"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").Value = "'0.0327"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +4.92%  "
$ws.Range("D51").Value = "'43.99"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -0.33%  "
